# Auto-generated edit script: update cryptos Price (D) and Volume(1h) (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '59.335.42'
$ws.Range("E2").Value = '  +0.47%  '
$ws.Range("D3").Value = '2.527.74'
$ws.Range("E3").Value = '  +1.26%  '
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("D5").Value = '''535.55'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.10%  '
$ws.Range("D6").Value = '''140.18'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.68%  '
$ws.Range("E7").Value = '  +0.42%  '
$ws.Range("E8").Value = '  -2.03%  '
$ws.Range("D9").Value = '2.535.04'
$ws.Range("E9").Value = '  +0.11%  '
$ws.Range("E10").Value = '  -0.67%  '
$ws.Range("E11").Value = '  +1.56%  '
$ws.Range("E12").Value = '  -1.45%  '
$ws.Range("E13").Value = '  +0.53%  '
$ws.Range("D14").Value = '2.976.65'
$ws.Range("E14").Value = '  +1.46%  '
$ws.Range("D15").Value = '''23.11'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -3.03%  '
$ws.Range("D16").Value = '59.284.35'
$ws.Range("E16").Value = '  +0.73%  '
$ws.Range("D17").Value = '''0.0000140'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.93%  '
$ws.Range("D18").Value = '2.555.04'
$ws.Range("E18").Value = '  +1.83%  '
$ws.Range("D19").Value = '''10.97'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.79%  '
$ws.Range("E20").Value = '  -1.51%  '
$ws.Range("D21").Value = '''321.13'
$ws.Range("D21").Style = "Normal"
$ws.Range("E22").Value = '  -0.05%  '
$ws.Range("D23").Value = '''5.80'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.78%  '
$ws.Range("D24").Value = '''62.47'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.24%  '
$ws.Range("D25").Value = '''0.420'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -4.50%  '
$ws.Range("E26").Value = '  +2.46%  '
$ws.Range("E27").Value = '  +0.58%  '
$ws.Range("D28").Value = '''7.85'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.15%  '
$ws.Range("D29").Value = '''6.73'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.68%  '
$ws.Range("D30").Value = '0.0₃0770'
$ws.Range("E30").Value = '  -1.00%  '
$ws.Range("E31").Value = '  +0.05%  '
$ws.Range("D32").Value = '''161.22'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.93%  '
$ws.Range("E33").Value = '  +0.42%  '
$ws.Range("D34").Value = '''1.13'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -5.75%  '
$ws.Range("E35").Value = '  -0.95%  '
$ws.Range("E36").Value = '  -0.43%  '
$ws.Range("E37").Value = '  -3.00%  '
$ws.Range("E38").Value = '  -2.86%  '
$ws.Range("D39").Value = '''37.07'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.59%  '
$ws.Range("E40").Value = '  -0.90%  '
$ws.Range("D41").Value = '''5.31'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -6.84%  '
$ws.Range("D42").Value = '''286.85'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -6.27%  '
$ws.Range("D43").Value = '''0.804'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.34%  '
$ws.Range("D44").Value = '''0.998'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.40%  '
$ws.Range("D45").Value = '''0.602'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.55%  '
$ws.Range("D46").Value = '''10.85'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.77%  '
$ws.Range("D47").Value = '''124.16'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.74%  '
$ws.Range("E48").Value = '  -0.54%  '
$ws.Range("D49").Value = '''18.60'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.28%  '
$ws.Range("D50").Value = '''0.0508'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.86%  '
$ws.Range("D51").Value = '''0.0222'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.52%  '
